# BOM.xlsx update: added mounting-hole hardware (Amphenol right-angle header)
# in place of the Samtec SMT header row, and bumped quantities for the
# mounting-hole related parts (rows 13, 14, 19) now that the board uses
# screw-mounted standoffs instead of sliding into the housing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Quantities that changed because mounting holes (and their hardware)
#     are now actually used ---
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 9
$ws.Range("D19").Value = 4

# --- Row 17 part swap: Samtec SMT header (SAM8980-ND) -> Amphenol
#     through-hole right-angle header (10129379-902004BLF-ND) ---
$ws.Range("B17").Value = "10129379-902004BLF-ND"
$ws.Range("C17").Value = "Connector Header Through Hole, Right Angle 2 position 0.100`" (2.54mm)"
$ws.Range("E17").Value = 0.1

$newUrl = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/10129379-902004BLF/7915980"
$ws.Range("H17").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("H17"), $newUrl) | Out-Null

# --- Selection left where the editor was last working ---
$ws.Activate() | Out-Null
$ws.Range("D20").Select() | Out-Null
